$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "./model_output/2025-08-20-19-55-06-None"
$ws.Range("B11").Value = 0.7050000000000001
